# Correct spacing used in ordered and unordered lists
#
# This script:
#  1. Removes the stray _GoBack bookmark that originally sat between
#     "Lists" and ":" in the opening paragraph.
#  2. Drops the trailing full stops from "Unordered sub-list.",
#     "And another item." and the final " item." run.
#  3. Removes the now-orphaned trailing space run at the very end of
#     the document and re-adds the _GoBack bookmark there instead,
#     right after "Final item".

$d = $word.ActiveDocument

# --- 1. Remove the original _GoBack bookmark -----------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Strip trailing periods from the two standalone list items --------
$d.Content.Find.Execute("Unordered sub-list.", $true, $false, $false, $false, $false, $true, 1, $false, "Unordered sub-list", 2) | Out-Null
$d.Content.Find.Execute("And another item.", $true, $false, $false, $false, $false, $true, 1, $false, "And another item", 2) | Out-Null

# --- 3. Fix the final paragraph ("Final item.") ---------------------------
# Locate "Final" so we only touch the period that belongs to the very last
# occurrence of " item." (there is another "... another item." earlier).
$rFinal = $d.Content
$rFinal.Find.Execute("Final", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rRest = $d.Range($rFinal.End, $d.Content.End)
$rRest.Find.Execute(" item.", $true, $false, $false, $false, $false, $true, 1, $false, " item", 2) | Out-Null

# Add the _GoBack bookmark immediately after "item" (collapsed range),
# placing it before the trailing space run is removed so that the
# insertion point is not the very last character of the document (doing
# it in that order avoids an edge-case with bookmarks collapsed at the
# absolute end of the document).
$full = $d.Content.Text
$pos = $full.Length - 2
$rBm = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $rBm) | Out-Null

# Now delete the orphaned trailing space run that used to follow "item."
$full2 = $d.Content.Text
$trailStart = $full2.Length - 2
$trailEnd = $full2.Length - 1
$rTrail = $d.Range($trailStart, $trailEnd)
if ($rTrail.Text -eq " ") {
    $rTrail.Delete() | Out-Null
}
